$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("level")

# --- Workbook window view change ---
$excel.ActiveWindow.Top = 4950

# --- Sheet view change: scroll so A43 is the top-left visible cell, select N52 ---
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("N52").Select()

# --- Cell value changes ---
    $ws.Range("I57").Value = 0.2
    $ws.Range("I58").Value = 0.2
    $ws.Range("I59").Value = 0.2
    $ws.Range("I60").Value = 0.2
    $ws.Range("I61").Value = 0.2
    $ws.Range("I62").Value = 0.2
    $ws.Range("I63").Value = 0.2
    $ws.Range("I64").Value = 0.2
    $ws.Range("I65").Value = 0.2
    $ws.Range("I66").Value = 0.2
    $ws.Range("I67").Value = 0.15
    $ws.Range("I68").Value = 0.15
    $ws.Range("I69").Value = 0.15
    $ws.Range("I70").Value = 0.15
    $ws.Range("I71").Value = 0.15
    $ws.Range("I72").Value = 0.15
    $ws.Range("I73").Value = 0.15
    $ws.Range("I74").Value = 0.15
    $ws.Range("I75").Value = 0.15
    $ws.Range("I76").Value = 0.15
    $ws.Range("I77").Value = 0.1
    $ws.Range("X77").Value = 2
    $ws.Range("I78").Value = 0.1
    $ws.Range("O78").Value = 46
    $ws.Range("X78").Value = 2
    $ws.Range("I79").Value = 0.1
    $ws.Range("O79").Value = 48
    $ws.Range("X79").Value = 2
    $ws.Range("I80").Value = 0.1
    $ws.Range("O80").Value = 50
    $ws.Range("X80").Value = 2
    $ws.Range("I81").Value = 0.1
    $ws.Range("O81").Value = 52
    $ws.Range("X81").Value = 2
    $ws.Range("I82").Value = 0.1
    $ws.Range("O82").Value = 54
    $ws.Range("X82").Value = 2
    $ws.Range("I83").Value = 0.1
    $ws.Range("O83").Value = 56
    $ws.Range("X83").Value = 2
    $ws.Range("I84").Value = 0.1
    $ws.Range("N84").Value = 37
    $ws.Range("O84").Value = 58
    $ws.Range("X84").Value = 2
    $ws.Range("I85").Value = 0.1
    $ws.Range("N85").Value = 37
    $ws.Range("O85").Value = 60
    $ws.Range("X85").Value = 2
    $ws.Range("I86").Value = 0.1
    $ws.Range("N86").Value = 37
    $ws.Range("O86").Value = 62
    $ws.Range("X86").Value = 2
    $ws.Range("I87").Value = 0.5
    $ws.Range("N87").Value = 39
    $ws.Range("O87").Value = 64
    $ws.Range("I88").Value = 0.5
    $ws.Range("N88").Value = 39
    $ws.Range("O88").Value = 66
    $ws.Range("I89").Value = 0.5
    $ws.Range("N89").Value = 39
    $ws.Range("O89").Value = 68
    $ws.Range("I90").Value = 0.5
    $ws.Range("N90").Value = 41
    $ws.Range("O90").Value = 70
    $ws.Range("I91").Value = 0.5
    $ws.Range("N91").Value = 41
    $ws.Range("O91").Value = 72
    $ws.Range("I92").Value = 0.5
    $ws.Range("N92").Value = 41
    $ws.Range("O92").Value = 74
    $ws.Range("I93").Value = 0.5
    $ws.Range("N93").Value = 43
    $ws.Range("O93").Value = 76
    $ws.Range("I94").Value = 0.5
    $ws.Range("N94").Value = 43
    $ws.Range("O94").Value = 78
    $ws.Range("I95").Value = 0.5
    $ws.Range("N95").Value = 43
    $ws.Range("O95").Value = 80
    $ws.Range("I96").Value = 0.5
    $ws.Range("N96").Value = 46
    $ws.Range("O96").Value = 82
    $ws.Range("I97").Value = 0.25
    $ws.Range("N97").Value = 46
    $ws.Range("O97").Value = 84
    $ws.Range("I98").Value = 0.25
    $ws.Range("N98").Value = 46
    $ws.Range("O98").Value = 86
    $ws.Range("I99").Value = 0.25
    $ws.Range("N99").Value = 51
    $ws.Range("O99").Value = 88
    $ws.Range("N100").Value = 51
    $ws.Range("O100").Value = 90
    $ws.Range("I101").Value = 0.15
    $ws.Range("N101").Value = 51
    $ws.Range("O101").Value = 92
    $ws.Range("I102").Value = 0.1
    $ws.Range("N102").Value = 55
    $ws.Range("O102").Value = 94

# --- Shared formula ref shrink (S16: Q16:S46 -> Q16:S30) ---
$ws.Range("S16").FormulaR1C1 = "=R[0]C[-1]-1"
